# Update the "Förändrad" (Changed) date column (column C) for every data
# row from row 2 through row 535: increment the stored date serial value
# from 45189 (2023-09-20) to 45190 (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = 535
$firstRow = 2

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
